$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (27) with the next schedule entry.
$ws.Range("A27").Value = "9/26"
$ws.Range("B27").Value = "11/21"
$ws.Range("C27").Value = "第77期 第二代星途"

# Match formatting of the preceding data rows (A/B text columns, C default).
$ws.Range("A27:B27").Style = $ws.Range("A26:B26").Style
$ws.Range("A27").NumberFormat = $ws.Range("A26").NumberFormat
$ws.Range("B27").NumberFormat = $ws.Range("B26").NumberFormat

# Update the active selection as in the saved workbook.
$ws.Range("C25").Select()
$excel.ActiveWindow.TopLeftCell = $ws.Range("A15")
